$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-18 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-19 Tuesday", 2)
$d.Content.Find.Execute("60×66=3960", $true, $false, $false, $false, $false, $true, 1, $false, "14×20=280", 2)
$d.Content.Find.Execute("73×81=5913", $true, $false, $false, $false, $false, $true, 1, $false, "28×51=1428", 2)
$d.Content.Find.Execute("99×83=8217", $true, $false, $false, $false, $false, $true, 1, $false, "15×99=1485", 2)
$d.Content.Find.Execute("19×28=532", $true, $false, $false, $false, $false, $true, 1, $false, "47×48=2256", 2)
$d.Content.Find.Execute("61×53=3233", $true, $false, $false, $false, $false, $true, 1, $false, "91×77=7007", 2)
$d.Content.Find.Execute("72×38=2736", $true, $false, $false, $false, $false, $true, 1, $false, "89×34=3026", 2)
$d.Content.Find.Execute("86×26=2236", $true, $false, $false, $false, $false, $true, 1, $false, "38×75=2850", 2)
$d.Content.Find.Execute("40×22=880", $true, $false, $false, $false, $false, $true, 1, $false, "92×34=3128", 2)
$d.Content.Find.Execute("48×16=768", $true, $false, $false, $false, $false, $true, 1, $false, "30×62=1860", 2)
$d.Content.Find.Execute("79×73=5767", $true, $false, $false, $false, $false, $true, 1, $false, "20×36=720", 2)
$d.Content.Find.Execute("98×38=3724", $true, $false, $false, $false, $false, $true, 1, $false, "39×57=2223", 2)
$d.Content.Find.Execute("81×57=4617", $true, $false, $false, $false, $false, $true, 1, $false, "20×64=1280", 2)
$d.Content.Find.Execute("43×99=4257", $true, $false, $false, $false, $false, $true, 1, $false, "56×78=4368", 2)
$d.Content.Find.Execute("42×28=1176", $true, $false, $false, $false, $false, $true, 1, $false, "57×16=912", 2)
$d.Content.Find.Execute("44×82=3608", $true, $false, $false, $false, $false, $true, 1, $false, "68×85=5780", 2)
$d.Content.Find.Execute("71×11=781", $true, $false, $false, $false, $false, $true, 1, $false, "38×53=2014", 2)
$d.Content.Find.Execute("79×33=2607", $true, $false, $false, $false, $false, $true, 1, $false, "93×24=2232", 2)
$d.Content.Find.Execute("62×39=2418", $true, $false, $false, $false, $false, $true, 1, $false, "93×89=8277", 2)
$d.Content.Find.Execute("73×80=5840", $true, $false, $false, $false, $false, $true, 1, $false, "25×61=1525", 2)
$d.Content.Find.Execute("39×90=3510", $true, $false, $false, $false, $false, $true, 1, $false, "57×90=5130", 2)
$d.Content.Find.Execute("51×50=2550", $true, $false, $false, $false, $false, $true, 1, $false, "86×45=3870", 2)
$d.Content.Find.Execute("89×79=7031", $true, $false, $false, $false, $false, $true, 1, $false, "48×58=2784", 2)
$d.Content.Find.Execute("39×91=3549", $true, $false, $false, $false, $false, $true, 1, $false, "87×28=2436", 2)
$d.Content.Find.Execute("41×96=3936", $true, $false, $false, $false, $false, $true, 1, $false, "20×46=920", 2)
$d.Content.Find.Execute("13×49=637", $true, $false, $false, $false, $false, $true, 1, $false, "14×12=168", 2)
